# ADD results from server
# Update row 2 values (A2:O2) on the "2025", "2030", and "2035" sheets with
# refreshed figures from the server.

$wb = $excel.ActiveWorkbook

# ---- Sheet "2025" ----
$ws = $wb.Worksheets.Item("2025")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0.003041639999999998
$ws.Range("E2").Value = 0.3192021974358237
$ws.Range("G2").Value = 0.2494892361374987
$ws.Range("I2").Value = 0.3425375850000001
$ws.Range("L2").Value = 0.5264838000000002
$ws.Range("M2").Value = 0.074078775
$ws.Range("N2").Value = 11.48276724557358
$ws.Range("O2").Value = 3.165850632917435

# ---- Sheet "2030" ----
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 0.01269103052680951
$ws.Range("B2").Value = 0.03782856521265217
$ws.Range("E2").Value = 0.1995122964266774
$ws.Range("I2").Value = 0.4358172464383578
$ws.Range("L2").Value = 0.0495198899999999
$ws.Range("M2").Value = 0.04272847500000002
$ws.Range("N2").Value = 4.791772732874055
$ws.Range("O2").Value = 2.043739180225848

# ---- Sheet "2035" ----
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 0.07897063915663666
$ws.Range("B2").Value = 0.02351531447219858
$ws.Range("E2").Value = 0.1529221346202301
$ws.Range("I2").Value = 0.4168164135616428
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0.04503457499999998
$ws.Range("N2").Value = 7.208983614552512
$ws.Range("O2").Value = 4.493642959409812
